$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 214, pushing the current
# row 214 (Asterix / 1a guarda lavada / 44432) down to row 216.
$ws.Rows.Item(214).Insert()
$ws.Rows.Item(214).Insert()

# New row 214: weekly update, Asterix, "1a (cosecha)"
$ws.Cells.Item(214, 1).Value2 = 11
$ws.Cells.Item(214, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(214, 3).Value2 = "Bíobío"
$ws.Cells.Item(214, 4).Value2 = 44656
$ws.Cells.Item(214, 5).Value2 = 8
$ws.Cells.Item(214, 6).Value2 = 100114001
$ws.Cells.Item(214, 7).Value2 = "Papa"
$ws.Cells.Item(214, 8).Value2 = "Asterix"
$ws.Cells.Item(214, 9).Value2 = "1a (cosecha)"
$ws.Cells.Item(214, 10).Value2 = 180
$ws.Cells.Item(214, 11).Value2 = 8000
$ws.Cells.Item(214, 12).Value2 = 8000
$ws.Cells.Item(214, 13).Value2 = 8000
$ws.Cells.Item(214, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(214, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(214, 16).Value2 = 320
$ws.Cells.Item(214, 17).Value2 = 25
$ws.Cells.Item(214, 18).Value2 = "Hortaliza"

# New row 215: weekly update, Rosara, "1a (cosecha)"
$ws.Cells.Item(215, 1).Value2 = 11
$ws.Cells.Item(215, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(215, 3).Value2 = "Bíobío"
$ws.Cells.Item(215, 4).Value2 = 44656
$ws.Cells.Item(215, 5).Value2 = 8
$ws.Cells.Item(215, 6).Value2 = 100114001
$ws.Cells.Item(215, 7).Value2 = "Papa"
$ws.Cells.Item(215, 8).Value2 = "Rosara"
$ws.Cells.Item(215, 9).Value2 = "1a (cosecha)"
$ws.Cells.Item(215, 10).Value2 = 180
$ws.Cells.Item(215, 11).Value2 = 7000
$ws.Cells.Item(215, 12).Value2 = 7000
$ws.Cells.Item(215, 13).Value2 = 7000
$ws.Cells.Item(215, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(215, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(215, 16).Value2 = 280
$ws.Cells.Item(215, 17).Value2 = 25
$ws.Cells.Item(215, 18).Value2 = "Hortaliza"
